$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "value" column was inserted right before the existing column B
# ("age"), pushing the old B:L columns over to C:M.
$ws.Range("B:B").Insert() | Out-Null

# Header for the freshly inserted column.
$ws.Range("B1").Value = "value"

# Data for the freshly inserted column.
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 30
$ws.Range("B4").Value = 40
$ws.Range("B5").Value = 5

# Leave the selection on the new column's header, matching the saved file.
$ws.Range("B1").Select() | Out-Null
